$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was collected for this market; insert it as the
# new top data row (row 2), pushing the existing rows down by one.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

$ws.Cells.Item(2, 1).Value = 5
$ws.Cells.Item(2, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(2, 3).Value = "Maule"
$ws.Cells.Item(2, 4).Value = 44462
$ws.Cells.Item(2, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2, 5).Value = 7
$ws.Cells.Item(2, 6).Value = 300000000
$ws.Cells.Item(2, 7).Value = "Espárragos"
$ws.Cells.Item(2, 8).Value = "Verde"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 2000
$ws.Cells.Item(2, 11).Value = 1800
$ws.Cells.Item(2, 12).Value = 2000
$ws.Cells.Item(2, 13).Value = 1900
$ws.Cells.Item(2, 14).Value = "`$/kilo"
$ws.Cells.Item(2, 15).Value = "Provincia de Linares"
$ws.Cells.Item(2, 16).Value = 1900
$ws.Cells.Item(2, 17).Value = 1
$ws.Cells.Item(2, 18).Value = "Hortaliza"
